# This workbook tracks daily "Ajo" (garlic) price records for
# Feria Lagunitas de Puerto Montt. The commit adds two new weekly
# price records at row 473, pushing the existing rows 473:573 down
# to 475:575 (dimension grows from A1:R573 to A1:R575).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 473; this shifts
# every column (A:R) of the existing rows 473:573 down to 475:575,
# and Excel grows the worksheet dimension/UsedRange automatically.
$ws.Rows("473:474").Insert()

# ---- New row 473 ----
$ws.Range("A473").Value2 = 4
$ws.Range("B473").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C473").Value2 = "Los Lagos"
$ws.Range("D473").Value2 = 45275
$ws.Range("E473").Value2 = 10
$ws.Range("F473").Value2 = 100112003
$ws.Range("G473").Value2 = "Ajo"
$ws.Range("H473").Value2 = "Chino"
$ws.Range("I473").Value2 = "Primera"
$ws.Range("J473").Value2 = 300
$ws.Range("K473").Value2 = 25000
$ws.Range("L473").Value2 = 26000
$ws.Range("M473").Value2 = 25500
$ws.Range("N473").Value2 = '$/caja 10 kilos'
$ws.Range("O473").Value2 = "China"
$ws.Range("P473").Value2 = 2550
$ws.Range("Q473").Value2 = 10
$ws.Range("R473").Value2 = "Hortaliza"

# ---- New row 474 ----
$ws.Range("A474").Value2 = 4
$ws.Range("B474").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C474").Value2 = "Los Lagos"
$ws.Range("D474").Value2 = 45275
$ws.Range("E474").Value2 = 10
$ws.Range("F474").Value2 = 100112003
$ws.Range("G474").Value2 = "Ajo"
$ws.Range("H474").Value2 = "Chino"
$ws.Range("I474").Value2 = "Primera"
$ws.Range("J474").Value2 = 100
$ws.Range("K474").Value2 = 28000
$ws.Range("L474").Value2 = 28000
$ws.Range("M474").Value2 = 28000
$ws.Range("N474").Value2 = '$/malla 10 kilos'
$ws.Range("O474").Value2 = "China"
$ws.Range("P474").Value2 = 2800
$ws.Range("Q474").Value2 = 10
$ws.Range("R474").Value2 = "Hortaliza"

# Apply the same date-formatted number format that column D uses
# elsewhere (style index 2 corresponds to the custom date/time format)
$ws.Range("D473").NumberFormat = $ws.Range("D472").NumberFormat
$ws.Range("D474").NumberFormat = $ws.Range("D472").NumberFormat
